$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-converted to a number by Excel
# are pre-formatted as Text so the stored value stays an exact string, matching
# the source data (which intentionally stores these dotted price strings as text).

$ws.Range("D5,D6,D7,D8,D9,D10,D11,D12,D15,D17,D20,D21,D22,D23,D25,D26,D27,D28,D29,D30,D32,D33,D35,D36,D37,D38,D40,D43,D47,D48,D50,D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.064.84"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.646.59"
$ws.Range("E3").Value = "  +0.14%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Value = "218.39"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6
$ws.Range("D6").Value = "0.5190"
$ws.Range("E6").Value = "  -0.19%  "

# Row 7
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").Value = "0.2626"
$ws.Range("E8").Value = "  +0.37%  "

# Row 9
$ws.Range("D9").Value = "0.06295"
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
$ws.Range("D10").Value = "20.24"
$ws.Range("E10").Value = "  -0.92%  "

# Row 11
$ws.Range("D11").Value = "0.07687"
$ws.Range("E11").Value = "  -0.94%  "

# Row 12
$ws.Range("D12").Value = "4.590"
$ws.Range("E12").Value = "  +2.67%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.655.20"
$ws.Range("E13").Value = "  +2.51%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.874.32"
$ws.Range("E14").Value = "  +0.20%  "

# Row 15
$ws.Range("D15").Value = "0.5566"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8098"
$ws.Range("E16").Value = "  +1.26%  "

# Row 17
$ws.Range("D17").Value = "65.07"
$ws.Range("E17").Value = "  +0.49%  "

# Row 18
$ws.Range("D18").Value = "26.060.99"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").Value = "4.601"
$ws.Range("E20").Value = "  -0.88%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "10.43"
$ws.Range("E21").Value = "  +3.23%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "192.56"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").Value = "5.907"
$ws.Range("E23").Value = "  -0.73%  "

# Row 24
$ws.Range("E24").Value = "  -0.21%  "

# Row 25
$ws.Range("D25").Value = "144.33"
$ws.Range("E25").Value = "  -1.27%  "

# Row 26
$ws.Range("D26").Value = "0.1180"
$ws.Range("E26").Value = "  -1.80%  "

# Row 27
$ws.Range("D27").Value = "7.175"
$ws.Range("E27").Value = "  +0.14%  "

# Row 28
$ws.Range("D28").Value = "15.82"
$ws.Range("E28").Value = "  -0.62%  "

# Row 29
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +1.71%  "

# Row 30
$ws.Range("D30").Value = "0.05336"
$ws.Range("E30").Value = "  -4.95%  "

# Row 31
$ws.Range("E31").Value = "  +0.41%  "

# Row 32
$ws.Range("D32").Value = "3.446"
$ws.Range("E32").Value = "  -0.16%  "

# Row 33
$ws.Range("D33").Value = "3.321"
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("E34").Value = "  -2.84%  "

# Row 35
$ws.Range("D35").Value = "2.417"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36
$ws.Range("D36").Value = "2.782"
$ws.Range("E36").Value = "  -0.39%  "

# Row 37
$ws.Range("D37").Value = "0.9403"
$ws.Range("E37").Value = "  +0.55%  "

# Row 38
$ws.Range("D38").Value = "0.5582"
$ws.Range("E38").Value = "  -1.60%  "

# Row 39
$ws.Range("E39").Value = "  -0.70%  "

# Row 40
$ws.Range("D40").Value = "5.767"
$ws.Range("E40").Value = "  -3.24%  "

# Row 41
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("D42").Value = "1.025.37"
$ws.Range("E42").Value = "  -2.38%  "

# Row 43
$ws.Range("D43").Value = "0.8251"
$ws.Range("E43").Value = "  -1.81%  "

# Row 44
$ws.Range("E44").Value = "  -1.45%  "

# Row 45
$ws.Range("D45").Value = "1.785.34"
$ws.Range("E45").Value = "  +0.18%  "

# Row 46
$ws.Range("E46").Value = "  +3.29%  "

# Row 47
$ws.Range("D47").Value = "57.25"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.75%  "

# Row 49
$ws.Range("E49").Value = "  -0.36%  "

# Row 50
$ws.Range("D50").Value = "7.889"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("D51").Value = "0.05098"
$ws.Range("E51").Value = "  -4.18%  "
